$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PAS-730: CA choice covered
# G3 previously held "MSRP_2000_CHOICE" -> now covers the new test case value.
$ws.Range("G3").Value = "MSRP_2000_CHOICE_TEST"

# Column G widened (best-fit) to accommodate the longer value.
$ws.Columns.Item(7).ColumnWidth = 22.28

# Cursor/selection left on E3 after the edit.
[void]$ws.Range("E3").Select()

# Page setup configured for printing (portrait orientation).
$ws.PageSetup.Orientation = 1
